# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds text-formatted numbers (e.g. "37.485.53"); force
# the cells to remain Text so Excel does not auto-convert them to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$updates = @(
    @{ Row = 2; Col = "D"; Value = "37.485.53" },
    @{ Row = 2; Col = "E"; Value = "  -0.75%  " },
    @{ Row = 3; Col = "D"; Value = "2.079.19" },
    @{ Row = 3; Col = "E"; Value = "  -0.01%  " },
    @{ Row = 4; Col = "E"; Value = "  +0.13%  " },
    @{ Row = 5; Col = "D"; Value = "232.69" },
    @{ Row = 5; Col = "E"; Value = "  -0.43%  " },
    @{ Row = 6; Col = "D"; Value = "0.633" },
    @{ Row = 6; Col = "E"; Value = "  +1.49%  " },
    @{ Row = 7; Col = "E"; Value = "  +0.00%  " },
    @{ Row = 8; Col = "D"; Value = "57.75" },
    @{ Row = 8; Col = "E"; Value = "  -1.35%  " },
    @{ Row = 9; Col = "D"; Value = "0.390" },
    @{ Row = 9; Col = "E"; Value = "  -0.97%  " },
    @{ Row = 10; Col = "E"; Value = "  -0.71%  " },
    @{ Row = 11; Col = "E"; Value = "  +2.50%  " },
    @{ Row = 12; Col = "D"; Value = "15.04" },
    @{ Row = 12; Col = "E"; Value = "  +1.71%  " },
    @{ Row = 13; Col = "D"; Value = "2.384.75" },
    @{ Row = 13; Col = "E"; Value = "  +0.01%  " },
    @{ Row = 14; Col = "D"; Value = "21.01" },
    @{ Row = 14; Col = "E"; Value = "  +0.50%  " },
    @{ Row = 15; Col = "E"; Value = "  -0.51%  " },
    @{ Row = 16; Col = "D"; Value = "5.32" },
    @{ Row = 16; Col = "E"; Value = "  -0.03%  " },
    @{ Row = 17; Col = "D"; Value = "2.073.28" },
    @{ Row = 17; Col = "E"; Value = "  -0.03%  " },
    @{ Row = 18; Col = "D"; Value = "37.450.39" },
    @{ Row = 18; Col = "E"; Value = "  -0.80%  " },
    @{ Row = 19; Col = "D"; Value = "70.66" },
    @{ Row = 19; Col = "E"; Value = "  -0.59%  " },
    @{ Row = 20; Col = "E"; Value = "  -2.28%  " },
    @{ Row = 21; Col = "D"; Value = "0.0₃0830" },
    @{ Row = 21; Col = "E"; Value = "  -0.45%  " },
    @{ Row = 22; Col = "D"; Value = "228.65" },
    @{ Row = 22; Col = "E"; Value = "  +0.09%  " },
    @{ Row = 23; Col = "E"; Value = "  -0.14%  " },
    @{ Row = 24; Col = "D"; Value = "2.36" },
    @{ Row = 24; Col = "E"; Value = "  -1.01%  " },
    @{ Row = 25; Col = "D"; Value = "2.37" },
    @{ Row = 25; Col = "E"; Value = "  -1.04%  " },
    @{ Row = 26; Col = "D"; Value = "9.68" },
    @{ Row = 26; Col = "E"; Value = "  +7.07%  " },
    @{ Row = 27; Col = "D"; Value = "170.11" },
    @{ Row = 27; Col = "E"; Value = "  -0.41%  " },
    @{ Row = 28; Col = "E"; Value = "  -4.10%  " },
    @{ Row = 29; Col = "D"; Value = "19.50" },
    @{ Row = 29; Col = "E"; Value = "  +0.13%  " },
    @{ Row = 30; Col = "E"; Value = "  -0.72%  " },
    @{ Row = 31; Col = "E"; Value = "  +0.78%  " },
    @{ Row = 32; Col = "D"; Value = "4.64" },
    @{ Row = 32; Col = "E"; Value = "  -1.13%  " },
    @{ Row = 33; Col = "E"; Value = "  +0.77%  " },
    @{ Row = 34; Col = "E"; Value = "  -0.41%  " },
    @{ Row = 35; Col = "D"; Value = "2.47" },
    @{ Row = 35; Col = "E"; Value = "  -0.32%  " },
    @{ Row = 36; Col = "D"; Value = "1.82" },
    @{ Row = 36; Col = "E"; Value = "  -0.12%  " },
    @{ Row = 37; Col = "D"; Value = "3.31" },
    @{ Row = 37; Col = "E"; Value = "  -2.54%  " },
    @{ Row = 38; Col = "E"; Value = "  +0.07%  " },
    @{ Row = 39; Col = "D"; Value = "5.31" },
    @{ Row = 39; Col = "E"; Value = "  +0.08%  " },
    @{ Row = 40; Col = "D"; Value = "0.0231" },
    @{ Row = 40; Col = "E"; Value = "  +7.50%  " },
    @{ Row = 41; Col = "D"; Value = "100.16" },
    @{ Row = 41; Col = "E"; Value = "  +0.59%  " },
    @{ Row = 44; Col = "D"; Value = "1.19" },
    @{ Row = 44; Col = "E"; Value = "  +3.04%  " },
    @{ Row = 45; Col = "D"; Value = "16.75" },
    @{ Row = 45; Col = "E"; Value = "  +2.04%  " },
    @{ Row = 46; Col = "D"; Value = "1.455.30" },
    @{ Row = 46; Col = "E"; Value = "  +0.23%  " },
    @{ Row = 47; Col = "E"; Value = "  -1.32%  " },
    @{ Row = 48; Col = "D"; Value = "3.97" },
    @{ Row = 48; Col = "E"; Value = "  -5.65%  " },
    @{ Row = 49; Col = "E"; Value = "  -2.08%  " },
    @{ Row = 50; Col = "D"; Value = "2.94" },
    @{ Row = 50; Col = "E"; Value = "  -2.24%  " },
    @{ Row = 51; Col = "D"; Value = "2.269.65" },
    @{ Row = 51; Col = "E"; Value = "  -0.01%  " }
)

foreach ($u in $updates) {
    $cellRef = "$($u.Col)$($u.Row)"
    $ws.Range($cellRef).Value = $u.Value
}

# Rows 42 and 43 swap coin identity (Cronos <-> HuobiToken) along with new price/volume data
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").Value = "2.91"
$ws.Range("E42").Value = "  +1.08%  "

$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "0.0953"
$ws.Range("E43").Value = "  -2.04%  "
